$p = $ppt.ActivePresentation

# --- 1. Table style changes on the three table slides -----------------
$newStyleId = "{657AC279-A813-4691-986C-15C855E3F17C}"
foreach ($slideIndex in 14,15,16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Swap the theme colour scheme (Integral -> Office) -------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB = 0   # dk1 000000
$tcs.Colors(2).RGB = 16777215   # lt1 FFFFFF
$tcs.Colors(3).RGB = 6968388   # dk2 44546A
$tcs.Colors(4).RGB = 15132391   # lt2 E7E6E6
$tcs.Colors(5).RGB = 13998939   # accent1 5B9BD5
$tcs.Colors(6).RGB = 3243501   # accent2 ED7D31
$tcs.Colors(7).RGB = 10855845   # accent3 A5A5A5
$tcs.Colors(8).RGB = 49407   # accent4 FFC000
$tcs.Colors(9).RGB = 12874308   # accent5 4472C4
$tcs.Colors(10).RGB = 4697456   # accent6 70AD47
$tcs.Colors(11).RGB = 12673797   # hlink 0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
